$wb = $excel.ActiveWorkbook

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1652.05
$ws.Range("I32").Value = 1437.6373
$ws.Range("J32").Value = 3820
$ws.Range("K32").Value = 1437.6373
$ws.Range("L32").Value = 3820
$ws.Range("M32").Value = -1150.6373
$ws.Range("N32").Value = -4394

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 16950918
$ws.Range("I63").Value = 1617.3478
$ws.Range("J63").Value = 76925370
$ws.Range("K63").Value = 1617.3478
$ws.Range("L63").Value = 76925370
$ws.Range("M63").Value = -931.3478
$ws.Range("N63").Value = -76926742

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 16950918
$ws.Range("I66").Value = 1617.3478
$ws.Range("J66").Value = 76925370
$ws.Range("K66").Value = 8086.739
$ws.Range("L66").Value = 384626850
$ws.Range("M66").Value = -4654.739
$ws.Range("N66").Value = -384633714

# BSM row 43
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 142342
$ws.Range("J43").Value = 142342
$ws.Range("L43").Value = 142342
$ws.Range("N43").Value = -142704

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 11905414
$ws.Range("J94").Value = 998.625
$ws.Range("L94").Value = 998.625
$ws.Range("N94").Value = -1900.625

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1462.5814
$ws.Range("J31").Value = 2005.8
$ws.Range("L31").Value = 2005.8
$ws.Range("N31").Value = -2595.8

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1462.5814
$ws.Range("J34").Value = 2005.8
$ws.Range("L34").Value = 2005.8
$ws.Range("N34").Value = -2409.8

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5716497.5
$ws.Range("I62").Value = 2306
$ws.Range("J62").Value = 20001976
$ws.Range("K62").Value = 2306
$ws.Range("L62").Value = 20001976
$ws.Range("M62").Value = -1682
$ws.Range("N62").Value = -20003224

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5716497.5
$ws.Range("I65").Value = 2306
$ws.Range("J65").Value = 20001976
$ws.Range("K65").Value = 11530
$ws.Range("L65").Value = 100009880
$ws.Range("M65").Value = -8410
$ws.Range("N65").Value = -100016120

# CRP row 95
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 20032.555
$ws.Range("J95").Value = 20032.555
$ws.Range("L95").Value = 20032.555
$ws.Range("N95").Value = -25524.555

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1291.9375
$ws.Range("I132").Value = 863.73914
$ws.Range("K132").Value = 2591.21742
$ws.Range("M132").Value = -61.21741999999995

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4182.1816
$ws.Range("J39").Value = 4267.1113
$ws.Range("L39").Value = 12801.3339
$ws.Range("N39").Value = -13389.3339

# CUL row 49
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 6004
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2092.625
$ws.Range("I55").Value = 477
$ws.Range("J55").Value = 2631.1667
$ws.Range("K55").Value = 1431
$ws.Range("L55").Value = 7893.500100000001
$ws.Range("M55").Value = -1254
$ws.Range("N55").Value = -8247.500100000001

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9280
$ws.Range("I80").Value = 20000
$ws.Range("J80").Value = 6600
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 6600
$ws.Range("M80").Value = -19002
$ws.Range("N80").Value = -8596

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 9280
$ws.Range("I83").Value = 20000
$ws.Range("J83").Value = 6600
$ws.Range("K83").Value = 100000
$ws.Range("L83").Value = 33000
$ws.Range("M83").Value = -95008
$ws.Range("N83").Value = -42984

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 635
$ws.Range("I97").Value = 621.36365
$ws.Range("K97").Value = 621.36365
$ws.Range("M97").Value = -125.36365

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 417.91666
$ws.Range("I107").Value = 314
$ws.Range("K107").Value = 314
$ws.Range("M107").Value = 1606

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2142.2954
$ws.Range("I132").Value = 1867.0312
$ws.Range("J132").Value = 2876.3333
$ws.Range("K132").Value = 5601.0936
$ws.Range("L132").Value = 8628.999899999999
$ws.Range("M132").Value = -3071.0936
$ws.Range("N132").Value = -13688.9999

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2212.7856
$ws.Range("I7").Value = 2138.3
$ws.Range("J7").Value = 2399
$ws.Range("K7").Value = 2138.3
$ws.Range("L7").Value = 2399
$ws.Range("M7").Value = -2026.3
$ws.Range("N7").Value = -2623

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1047.625
$ws.Range("I22").Value = 890.5
$ws.Range("K22").Value = 890.5
$ws.Range("M22").Value = -595.5

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1047.625
$ws.Range("I27").Value = 890.5
$ws.Range("K27").Value = 890.5
$ws.Range("M27").Value = -783.5

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3969
$ws.Range("I40").Value = 2534.75
$ws.Range("J40").Value = 6837.5
$ws.Range("K40").Value = 2534.75
$ws.Range("L40").Value = 6837.5
$ws.Range("M40").Value = -2398.75
$ws.Range("N40").Value = -7109.5

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4056.8333
$ws.Range("I46").Value = 694
$ws.Range("K46").Value = 694
$ws.Range("M46").Value = -506

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 217.32
$ws.Range("I55").Value = 111.5
$ws.Range("J55").Value = 993.3333
$ws.Range("K55").Value = 111.5
$ws.Range("L55").Value = 993.3333
$ws.Range("M55").Value = 61.5
$ws.Range("N55").Value = -1339.3333

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1589.2
$ws.Range("J82").Value = 1623.1428
$ws.Range("L82").Value = 1623.1428
$ws.Range("N82").Value = -2345.1428

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1589.2
$ws.Range("J85").Value = 1623.1428
$ws.Range("L85").Value = 1623.1428
$ws.Range("N85").Value = -4119.1428

# LTW row 87
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 21999.5
$ws.Range("J87").Value = 21999.5
$ws.Range("L87").Value = 21999.5
$ws.Range("N87").Value = -24245.5

# LTW row 90
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H90").Value = 21999.5
$ws.Range("J90").Value = 21999.5
$ws.Range("L90").Value = 65998.5
$ws.Range("N90").Value = -77230.5

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1283.3334
$ws.Range("I100").Value = 1133.3334
$ws.Range("J100").Value = 1433.3334
$ws.Range("K100").Value = 1133.3334
$ws.Range("L100").Value = 1433.3334
$ws.Range("M100").Value = -592.3334
$ws.Range("N100").Value = -2515.3334

# LTW row 109
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 41668800
$ws.Range("I122").Value = 62502000
$ws.Range("J122").Value = 2402.5
$ws.Range("K122").Value = 187506000
$ws.Range("L122").Value = 7207.5
$ws.Range("M122").Value = -187503550
$ws.Range("N122").Value = -12107.5

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2212.7856
$ws.Range("I126").Value = 2138.3
$ws.Range("J126").Value = 2399
$ws.Range("K126").Value = 6414.900000000001
$ws.Range("L126").Value = 7197
$ws.Range("M126").Value = -3944.900000000001
$ws.Range("N126").Value = -12137

# WVR row 51
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 12000
$ws.Range("J51").Value = 12000
$ws.Range("L51").Value = 12000
$ws.Range("N51").Value = -13020

# WVR row 57
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 35000
$ws.Range("J57").Value = 35000
$ws.Range("L57").Value = 35000
$ws.Range("N57").Value = -36508

# WVR row 92
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 18942.428
$ws.Range("J92").Value = 18942.428
$ws.Range("L92").Value = 18942.428
$ws.Range("N92").Value = -23934.428

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1274.8096
$ws.Range("I96").Value = 1371.3334
$ws.Range("J96").Value = 1033.5
$ws.Range("K96").Value = 1371.3334
$ws.Range("L96").Value = 1033.5
$ws.Range("M96").Value = 1.666600000000017
$ws.Range("N96").Value = -3779.5

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 351.75
$ws.Range("I107").Value = 285.75
$ws.Range("K107").Value = 857.25
$ws.Range("M107").Value = 1062.75

# WVR row 109
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 38569.2
$ws.Range("J109").Value = 35626
$ws.Range("L109").Value = 35626
$ws.Range("N109").Value = -38400
